$wb = $excel.ActiveWorkbook

# --- Content updates -------------------------------------------------
# The "Status" text (shared across the Overview sheet and both language
# sheets) moves from "Ready for handoff" to the new handback status now
# that the report has been generated post-handback.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Range("K2").Value = "2016-08-22 18:52:49"
$zhcn.Range("P2").Value = ""

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Range("K2").Value = "2016-08-22 18:52:57"
$dede.Range("P2").Value = ""

# --- Column width refresh --------------------------------------------
# After the text changes above, the "Status" and "Error Detail" columns
# are best-fit to their new contents (longer status text widens the
# Status columns, the now-empty error detail shrinks that column back to
# roughly the width of its header).
$overview.Range("E:F").ColumnWidth = 29.166666666666668
$zhcn.Range("C:C").ColumnWidth = 29.166666666666668
$zhcn.Range("P:P").ColumnWidth = 12.833333333333334
$dede.Range("C:C").ColumnWidth = 29.166666666666668
$dede.Range("P:P").ColumnWidth = 12.833333333333334
